$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold numeric-looking text (e.g. "0.2860", "1.001")
# that must stay literal text. If assigned while the cell is General-
# formatted, Excel (like real Excel COM) reinterprets the string as a
# number and normalises/truncates it (e.g. "0.2860" -> 0.286). Force
# the cell to Text first, assign, then restore the default "Normal"
# style so no stray number-format is left behind.
$textForceCells = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12",
    "D14", "D15", "D16", "D20", "D22", "D23", "D24", "D26",
    "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35",
    "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43",
    "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.249.80"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.862.74"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "236.88"
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").Value = "0.4677"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "0.2860"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").Value = "0.06534"
$ws.Range("E9").Value = "  -0.30%  "

$ws.Range("D10").Value = "22.01"
$ws.Range("E10").Value = "  +8.98%  "

$ws.Range("D11").Value = "0.07918"
$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("D12").Value = "97.53"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("D13").Value = "1.869.32"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").Value = "5.165"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").Value = "0.6824"
$ws.Range("E15").Value = "  +1.81%  "

$ws.Range("D16").Value = "270.74"
$ws.Range("E16").Value = "  -3.38%  "

$ws.Range("D17").Value = "30.246.84"

$ws.Range("E18").Value = "  +7.04%  "

$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").Value = "0.000007354"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("D21").Value = "2.111.81"

$ws.Range("D22").Value = "5.325"
$ws.Range("E22").Value = "  -2.37%  "

$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").Value = "6.184"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").Value = "9.229"
$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("D27").Value = "18.93"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").Value = "1.954"
$ws.Range("E28").Value = "  +2.52%  "

$ws.Range("E29").Value = "  +3.05%  "

$ws.Range("D30").Value = "0.09837"
$ws.Range("E30").Value = "  +2.38%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.363"
$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.482"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").Value = "4.057"
$ws.Range("E33").Value = "  -1.51%  "

$ws.Range("D34").Value = "0.04714"
$ws.Range("E34").Value = "  +1.12%  "

$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  +3.28%  "

$ws.Range("D36").Value = "0.7020"
$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").Value = "0.01878"
$ws.Range("E38").Value = "  +1.26%  "

$ws.Range("D39").Value = "2.626"
$ws.Range("E39").Value = "  +3.82%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "75.57"
$ws.Range("E40").Value = "  +4.26%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.277"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").Value = "1.946"
$ws.Range("E42").Value = "  +0.90%  "

$ws.Range("D43").Value = "0.8521"

$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").Value = "0.4164"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "103.20"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "7.183"
$ws.Range("E47").Value = "  +0.62%  "

$ws.Range("D48").Value = "952.09"
$ws.Range("E48").Value = "  -3.64%  "

$ws.Range("D49").Value = "9.220"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").Value = "34.18"
$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").Value = "0.05652"
$ws.Range("E51").Value = "  +0.29%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
